$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.502379536628723
$ws.Range("B1").Value = 2.776848793029785
$ws.Range("C1").Value = 3.043171644210815
$ws.Range("D1").Value = 3.066450595855713
$ws.Range("E1").Value = 2.397652626037598
